$d = $word.ActiveDocument

# Collapse to the very end of the document body content.
$end = $d.Content
$end.Collapse(0)  # wdCollapseEnd

# Insert the new "SS manual V2" paragraph (as two runs: "SS manual V" + "2")
# followed by a new empty paragraph, via a flat-OPC WordOpenXML fragment so
# the two text runs land as separate <w:r> elements and the trailing
# paragraph is truly empty (<w:p/>), matching how Word represents it.
$xml = '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:r><w:t>SS manual V</w:t></w:r><w:r><w:t>2</w:t></w:r></w:p><w:p/></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$end.InsertXML($xml)
